# progress_tracking.xlsx — "adding more data for Q&A"
#
# 1) Remove the now-unused "Daily Problems" sheet.
# 2) On "Introduction to Python", add a new column F with header
#    "#Don't touch Medium Questions yet" (bold, matching the other headers).
# 3) Update a couple of Remarks cells (E4, E6) from "Try Again" to "Good".
# 4) Leave the cursor on A26 (matches the author's last selection).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Introduction to Python")

# --- remove the empty "Daily Problems" sheet -------------------------------
$wb.Worksheets.Item("Daily Problems").Delete()

# --- new column F header, styled like the rest of row 1 --------------------
$ws.Range("F1").Value = "#Don't touch Medium Questions yet"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Columns.Item(6).ColumnWidth = 30.830729166666668   # ~31.6640625 raw width

# --- mark two more attempts as solved ("Try Again" -> "Good") --------------
$ws.Range("E4").Value = "Good"
$ws.Range("E6").Value = "Good"

# --- restore the author's selection -----------------------------------------
$ws.Range("A26").Select()
